$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in column C (header is in row 1)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

# Column C ("Förändrad") holds a date serial value of 45190 for every
# data row; update it to 45192 for all rows that currently have 45190.
$rng = $ws.Range("C2:C$lastRow")
foreach ($cell in $rng.Cells) {
    if ($cell.Value2 -eq 45190) {
        $cell.Value2 = 45192
    }
}
